$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D2 and D3 to be quote-prefixed text values ("5" and "6")
$ws.Range("D2").Value = "'5"
$ws.Range("D3").Value = "'6"

# Update the active selection to D4
$ws.Range("D4").Select()
